# Scheduled runner update: refresh cached market-board figures
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns H-N)
# across the per-class Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2462.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2462.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 7387.5
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -7723.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1887.3889
$ws.Range("I28").Value = 783
$ws.Range("K28").Value = 783
$ws.Range("M28").Value = -298

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 210.25
$ws.Range("I33").Value = 210.25
$ws.Range("K33").Value = 210.25
$ws.Range("M33").Value = 18.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 10954
$ws.Range("I137").Value = 1998.625
$ws.Range("J137").Value = 34835
$ws.Range("K137").Value = 5995.875
$ws.Range("L137").Value = 104505
$ws.Range("M137").Value = -3445.875
$ws.Range("N137").Value = -109605

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3686.948
$ws.Range("I32").Value = 2998.5754
$ws.Range("J32").Value = 16249.75
$ws.Range("K32").Value = 2998.5754
$ws.Range("L32").Value = 16249.75
$ws.Range("M32").Value = -2711.5754
$ws.Range("N32").Value = -16823.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 41245.363
$ws.Range("I45").Value = 54026.125
$ws.Range("J45").Value = 7163.3335
$ws.Range("K45").Value = 54026.125
$ws.Range("L45").Value = 7163.3335
$ws.Range("M45").Value = -53649.125
$ws.Range("N45").Value = -7917.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4491.7075
$ws.Range("I122").Value = 4264.8486
$ws.Range("J122").Value = 5427.5
$ws.Range("K122").Value = 12794.5458
$ws.Range("L122").Value = 16282.5
$ws.Range("M122").Value = -10344.5458
$ws.Range("N122").Value = -21182.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2579.1785
$ws.Range("I132").Value = 1878.1666
$ws.Range("J132").Value = 3841
$ws.Range("K132").Value = 5634.4998
$ws.Range("L132").Value = 11523
$ws.Range("M132").Value = -3104.4998
$ws.Range("N132").Value = -16583

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 993.5455
$ws.Range("I22").Value = 1012.7
$ws.Range("J22").Value = 802
$ws.Range("K22").Value = 1012.7
$ws.Range("L22").Value = 802
$ws.Range("M22").Value = -839.7
$ws.Range("N22").Value = -1148

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3693.7334
$ws.Range("I86").Value = 3504.7273
$ws.Range("K86").Value = 3504.7273
$ws.Range("M86").Value = -2381.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3693.7334
$ws.Range("I89").Value = 3504.7273
$ws.Range("K89").Value = 17523.6365
$ws.Range("M89").Value = -11907.6365

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 62500652
$ws.Range("I94").Value = 62500652
$ws.Range("K94").Value = 62500652
$ws.Range("M94").Value = -62500201

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 102034.8
$ws.Range("I99").Value = 143471.14
$ws.Range("K99").Value = 143471.14
$ws.Range("M99").Value = -141973.14

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2910.5557
$ws.Range("I105").Value = 1500
$ws.Range("K105").Value = 1500
$ws.Range("M105").Value = 247

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 14289610
$ws.Range("I132").Value = 18522034
$ws.Range("J132").Value = 5176.625
$ws.Range("K132").Value = 55566102
$ws.Range("L132").Value = 15529.875
$ws.Range("M132").Value = -55563572
$ws.Range("N132").Value = -20589.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3010.7778
$ws.Range("I5").Value = 601.7143
$ws.Range("K5").Value = 1805.1429
$ws.Range("M5").Value = -1693.1429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2866.7646
$ws.Range("I132").Value = 1675
$ws.Range("J132").Value = 3516.818
$ws.Range("K132").Value = 15075
$ws.Range("L132").Value = 31651.362
$ws.Range("M132").Value = -12545
$ws.Range("N132").Value = -36711.362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 3010.7778
$ws.Range("I135").Value = 601.7143
$ws.Range("K135").Value = 5415.428699999999
$ws.Range("M135").Value = -2880.428699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 142859460
$ws.Range("I80").Value = 166668290
$ws.Range("K80").Value = 166668290
$ws.Range("M80").Value = -166667292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 142859460
$ws.Range("I83").Value = 166668290
$ws.Range("K83").Value = 833341450
$ws.Range("M83").Value = -833336458

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5500.143
$ws.Range("I113").Value = 5617.1763
$ws.Range("J113").Value = 5002.75
$ws.Range("K113").Value = 5617.1763
$ws.Range("L113").Value = 5002.75
$ws.Range("M113").Value = -3447.1763
$ws.Range("N113").Value = -9342.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 15394614
$ws.Range("I122").Value = 76923070
$ws.Range("K122").Value = 230769210
$ws.Range("M122").Value = -230766760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2386.7058
$ws.Range("I132").Value = 1664.3077
$ws.Range("J132").Value = 4734.5
$ws.Range("K132").Value = 4992.9231
$ws.Range("L132").Value = 14203.5
$ws.Range("M132").Value = -2462.9231
$ws.Range("N132").Value = -19263.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 10350.333
$ws.Range("J136").Value = 10350.333
$ws.Range("L136").Value = 31050.999
$ws.Range("N136").Value = -36150.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 835.1429
$ws.Range("I16").Value = 872.53845
$ws.Range("J16").Value = 349
$ws.Range("K16").Value = 872.53845
$ws.Range("L16").Value = 349
$ws.Range("M16").Value = -702.53845
$ws.Range("N16").Value = -689

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 637.375
$ws.Range("I22").Value = 596.2
$ws.Range("K22").Value = 596.2
$ws.Range("M22").Value = -301.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 637.375
$ws.Range("I27").Value = 596.2
$ws.Range("K27").Value = 596.2
$ws.Range("M27").Value = -489.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1405.12
$ws.Range("I61").Value = 1221.9
$ws.Range("K61").Value = 1221.9
$ws.Range("M61").Value = -1019.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1747.4166
$ws.Range("I93").Value = 1815.3636
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 1815.3636
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -567.3635999999999
$ws.Range("N93").Value = -3496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1084691
$ws.Range("I100").Value = 1290522.8
$ws.Range("K100").Value = 1290522.8
$ws.Range("M100").Value = -1289981.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1405.12
$ws.Range("I113").Value = 1221.9
$ws.Range("K113").Value = 1221.9
$ws.Range("M113").Value = 948.0999999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 25837.666
$ws.Range("J45").Value = 25837.666
$ws.Range("L45").Value = 25837.666
$ws.Range("N45").Value = -26819.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 27450
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 27450
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 19500
$ws.Range("J99").Value = 19500
$ws.Range("L99").Value = 19500
$ws.Range("N99").Value = -25490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 41669524
$ws.Range("I122").Value = 3049.5
$ws.Range("K122").Value = 9148.5
$ws.Range("M122").Value = -6698.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 33340164
$ws.Range("I132").Value = 41674144
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 125022432
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -125019902
$ws.Range("N132").Value = -17810
